# ajuste: corrigindo as categorias
#
# Adds a "Total" column (T) with the row sum for each existing data row,
# and adds two new rows:
#   row 7 -> "Outros" category with its per-age-group counts + row total
#   row 8 -> "Total" category (grand total) with per-age-group sums + grand total
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the "Total" column
$ws.Range("T1").Value = "Total"

# Row totals for the existing categories (rows 2-6)
$ws.Range("T2").Value = 81182
$ws.Range("T3").Value = 7270
$ws.Range("T4").Value = 34679
$ws.Range("T5").Value = 13535
$ws.Range("T6").Value = 47276

# New row 7: "Outros"
$ws.Range("A7").Value = "Outros"
$ws.Range("B7").Value = 7237
$ws.Range("C7").Value = 395
$ws.Range("D7").Value = 570
$ws.Range("E7").Value = 2239
$ws.Range("F7").Value = 3466
$ws.Range("G7").Value = 3692
$ws.Range("H7").Value = 3813
$ws.Range("I7").Value = 4028
$ws.Range("J7").Value = 4609
$ws.Range("K7").Value = 5120
$ws.Range("L7").Value = 5493
$ws.Range("M7").Value = 5435
$ws.Range("N7").Value = 5185
$ws.Range("O7").Value = 5038
$ws.Range("P7").Value = 5309
$ws.Range("Q7").Value = 6250
$ws.Range("R7").Value = 17849
$ws.Range("S7").Value = 697
$ws.Range("T7").Value = 86425

# New row 8: "Total" (grand total row)
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = 8199
$ws.Range("C8").Value = 618
$ws.Range("D8").Value = 857
$ws.Range("E8").Value = 2753
$ws.Range("F8").Value = 4183
$ws.Range("G8").Value = 4792
$ws.Range("H8").Value = 5507
$ws.Range("I8").Value = 6609
$ws.Range("J8").Value = 9067
$ws.Range("K8").Value = 12236
$ws.Range("L8").Value = 16399
$ws.Range("M8").Value = 19615
$ws.Range("N8").Value = 22130
$ws.Range("O8").Value = 23482
$ws.Range("P8").Value = 26951
$ws.Range("Q8").Value = 30162
$ws.Range("R8").Value = 75907
$ws.Range("S8").Value = 900
$ws.Range("T8").Value = 270367
